# Add two new result sheets (NC11, HC6) at the end of the workbook, mirroring
# the existing NC10/HC5 "summary" sheets (same B1:E1 headers, A2 label, and
# B2:E2 numeric values laid out the same way).
#
# NC11 is an exact duplicate of NC10 (same values), and HC6 duplicates HC5's
# layout/formatting but carries new numbers for the "Holding control" run.

$wb = $excel.ActiveWorkbook

# Remember the sheet that's active right now so we can restore the selection
# afterwards (adding/copying sheets moves the active tab).
$originalActive = $wb.ActiveSheet

# --- NC11: duplicate of NC10 -------------------------------------------------
$ncSource = $wb.Worksheets.Item("NC10")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ncSource.Copy($null, $lastSheet)
$nc11 = $wb.Worksheets.Item($wb.Worksheets.Count)
$nc11.Name = "NC11"

# --- HC6: duplicate of HC5, with updated values ------------------------------
$hcSource = $wb.Worksheets.Item("HC5")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$hcSource.Copy($null, $lastSheet)
$hc6 = $wb.Worksheets.Item($wb.Worksheets.Count)
$hc6.Name = "HC6"

$hc6.Range("B2").Value = 2417.211933973449
$hc6.Range("C2").Value = 12318.81187969196
$hc6.Range("D2").Value = 133.6515102396718
$hc6.Range("E2").Value = 14869.67532390508

# Restore original selection/active sheet.
$originalActive.Activate()
